$wb = $excel.ActiveWorkbook

# --- Sheet "Summary" ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.2330960854092527
$ws1.Range("C2").Value = 0.06100217864923747
$ws1.Range("D2").Value = 1
$ws1.Range("E2").Value = 0.1149897330595483
$ws1.Range("F2").Value = 0.245183887915937
$ws1.Range("G2").Value = 0.6281276962899051
$ws1.Range("H2").Value = 0.7077314071696094
$ws1.Range("I2").Value = 28
$ws1.Range("J2").Value = 431
$ws1.Range("K2").Value = 103
$ws1.Range("L2").Value = 0

# --- Sheet "Classification Report" ---
$ws2 = $wb.Worksheets.Item("Classification Report")

# Row 2 ("0")
$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = 0.1928838951310861
$ws2.Range("D2").Value = 0.3233908948194663

# Row 3 ("1")
$ws2.Range("B3").Value = 0.06100217864923747
$ws2.Range("C3").Value = 1
$ws2.Range("D3").Value = 0.1149897330595483

# Row 4 ("accuracy")
$ws2.Range("B4").Value = 0.2330960854092527
$ws2.Range("C4").Value = 0.2330960854092527
$ws2.Range("D4").Value = 0.2330960854092527
$ws2.Range("E4").Value = 0.2330960854092527

# Row 5 ("macro avg")
$ws2.Range("B5").Value = 0.5305010893246187
$ws2.Range("C5").Value = 0.596441947565543
$ws2.Range("D5").Value = 0.2191903139395072

# Row 6 ("weighted avg")
$ws2.Range("B6").Value = 0.9532171903953356
$ws2.Range("C6").Value = 0.2330960854092527
$ws2.Range("D6").Value = 0.313007918788723

# --- Sheet "Confusion Matrix" ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

# Row 2 ("Actual 0")
$ws3.Range("B2").Value = 103
$ws3.Range("C2").Value = 431

# Row 3 ("Actual 1")
$ws3.Range("B3").Value = 0
$ws3.Range("C3").Value = 28
